# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across achievement / work
# experience bullet paragraphs, matching the target diff exactly.

$d = $word.ActiveDocument

# Color value for RGB(0x2C, 0x3E, 0x50) packed as Word's BGR-order Long.
$metricColor = 5258796

function Highlight-Metrics($Paragraph, $Terms) {
    foreach ($term in $Terms) {
        $rng = $Paragraph.Range
        $found = $rng.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = $true
            $rng.Font.Color = $metricColor
        }
    }
}

# Map of a distinctive, unique substring for locating each target paragraph
# to the ordered list of sub-strings within it that must become bold+colored.
$plan = @(
    @{ Match = "Discovered systematic race coding errors"; Terms = @("23%", "64%") },
    @{ Match = [char]0x2022 + " Utilized advanced sampling methods to decrease survey margin of error from "; Terms = @([char]0xB1 + "4.2%", [char]0xB1 + "2.1%", "71%", "87%") },
    @{ Match = "Trigonometric algorithm for boundary estimation"; Terms = @("73.5%", "$" + "4.7M") },
    @{ Match = "Built real-time FEC analysis systems"; Terms = @("$" + "2") },
    @{ Match = "Modernized legacy ETL processes"; Terms = @("57%") },
    @{ Match = "178% accuracy improvement"; Terms = @("178%") },
    @{ Match = "Algorithmic innovation: Pioneered trigonometric boundary estimation"; Terms = @("73.5%") },
    @{ Match = "$" + "4.7M savings enabled nonprofit access"; Terms = @("$" + "4.7M") },
    @{ Match = "Platform impact: Built redistricting system serving"; Terms = @("12,847") },
    @{ Match = "Predictive excellence: Utilized advanced sampling methods"; Terms = @([char]0xB1 + "4.2%", [char]0xB1 + "2.1%") },
    @{ Match = "Increased voter turnout prediction accuracy from"; Terms = @("71%", "87%") }
)

foreach ($item in $plan) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($item.Match)) {
            Highlight-Metrics $p $item.Terms
            break
        }
    }
}
